# repull data, push all data, mean calculation
# Update the dSF column (F) with newly re-pulled values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = -4
    3  = -1
    4  = 4
    5  = -2
    6  = 0
    7  = -5
    8  = 1
    9  = -8
    10 = 4
    11 = -6
    12 = -3
    13 = -2
    14 = 3
    15 = -5
    17 = -4
    18 = -3
    19 = -2
    20 = 3
    21 = -1
    22 = 1
    23 = -5
    24 = 3
    25 = -1
    26 = -1
    27 = -1
    28 = 3
    29 = -1
    30 = 1
    31 = -1
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 6).Value = $newValues[$row]
}
